# Adds a new "2022-Q4" worksheet (right after "总计") with fresh fund-holding
# data, and updates the "总计" (summary) sheet so its top row reflects the
# new quarter while the previously-existing rows shift down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Locate the existing "总计" summary sheet and the "2022-Q3" sheet
#    (used purely as a formatting template for the new sheet).
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Item("2022-Q3")

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q4" sheet right after "总计" (i.e. right before
#    the sheet that is currently "2022-Q3"). This matches the workbook.xml
#    ordering: 总计, 2022-Q4, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q1, 2020-Q4.
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $summary)
$q4.Name = "2022-Q4"

# Match the page-margin conventions used by the sibling quarterly sheets
# (PageSetup margins are expressed in points; 72pt = 1in).
$q4.PageSetup.LeftMargin = 54
$q4.PageSetup.RightMargin = 54
$q4.PageSetup.TopMargin = 72
$q4.PageSetup.BottomMargin = 72
$q4.PageSetup.HeaderMargin = 36
$q4.PageSetup.FooterMargin = 36

# Header row (column B..H) -- bold + bordered + centered, mirroring the
# look of the other quarterly sheets.
$headerVals = New-Object 'object[,]' 1,7
$headerVals[0,0] = "基金代码"
$headerVals[0,1] = "基金名称"
$headerVals[0,2] = "基金规模"
$headerVals[0,3] = "股票总仓位"
$headerVals[0,4] = "仓位占比"
$headerVals[0,5] = "持有市值(亿元)"
$headerVals[0,6] = "仓位排名"

$headerRng = $q4.Range("B1:H1")
$headerRng.Value = $headerVals
$headerRng.Font.Bold = $true
$headerRng.Borders.LineStyle = 1
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4160

# Column A index values (0,1,2) -- same bold/bordered/centered look as the
# header, matching the other sheets' "A" column styling.
$aRng = $q4.Range("A2:A4")
$aVals = New-Object 'object[,]' 3,1
$aVals[0,0] = 0
$aVals[1,0] = 1
$aVals[2,0] = 2
$aRng.Value = $aVals
$aRng.Font.Bold = $true
$aRng.Borders.LineStyle = 1
$aRng.HorizontalAlignment = -4108
$aRng.VerticalAlignment = -4160

# Data rows 2-4 (B..G are stored as text in every quarterly sheet, even
# though they look numeric -- e.g. fund codes with leading zeros, and
# percentages kept as plain text strings). Column H (rank) is numeric.
$bgRng = $q4.Range("B2:G4")
$bgVals = New-Object 'object[,]' 3,6
$bgVals[0,0] = "'162102"; $bgVals[0,1] = "金鹰中小盘精选混合"; $bgVals[0,2] = "'3.48"; $bgVals[0,3] = "'76.23"; $bgVals[0,4] = "'5.64"; $bgVals[0,5] = "'0.1963"
$bgVals[1,0] = "'210009"; $bgVals[1,1] = "金鹰核心资源混合";   $bgVals[1,2] = "'3.14"; $bgVals[1,3] = "'91.78"; $bgVals[1,4] = "'6.25"; $bgVals[1,5] = "'0.1962"
$bgVals[2,0] = "'001167"; $bgVals[2,1] = "金鹰科技创新股票";   $bgVals[2,2] = "'3.17"; $bgVals[2,3] = "'91.02"; $bgVals[2,4] = "'5.84"; $bgVals[2,5] = "'0.1851"
$bgRng.Value = $bgVals

$hRng = $q4.Range("H2:H4")
$hVals = New-Object 'object[,]' 3,1
$hVals[0,0] = 2
$hVals[1,0] = 3
$hVals[2,0] = 3
$hRng.Value = $hVals

$q4.Range("A1").Select()

# ---------------------------------------------------------------------
# 3) Update the "总计" sheet: insert the new 2022-Q4 row at the top of the
#    data (row 2) and push the previously existing rows down by one.
# ---------------------------------------------------------------------
$summaryVals = New-Object 'object[,]' 6,3
$summaryVals[0,0] = "2022-Q4"; $summaryVals[0,1] = 3; $summaryVals[0,2] = 0.58
$summaryVals[1,0] = "2022-Q3"; $summaryVals[1,1] = 3; $summaryVals[1,2] = 0.43
$summaryVals[2,0] = "2022-Q2"; $summaryVals[2,1] = 3; $summaryVals[2,2] = 0.4
$summaryVals[3,0] = "2022-Q1"; $summaryVals[3,1] = 2; $summaryVals[3,2] = 0.34
$summaryVals[4,0] = "2021-Q1"; $summaryVals[4,1] = 3; $summaryVals[4,2] = 0.07000000000000001
$summaryVals[5,0] = "2020-Q4"; $summaryVals[5,1] = 2; $summaryVals[5,2] = 0.02

$summary.Range("B2:D7").Value = $summaryVals

$aIdxRng = $summary.Range("A2:A7")
$aIdxVals = New-Object 'object[,]' 6,1
$aIdxVals[0,0] = 0
$aIdxVals[1,0] = 1
$aIdxVals[2,0] = 2
$aIdxVals[3,0] = 3
$aIdxVals[4,0] = 4
$aIdxVals[5,0] = 5
$aIdxRng.Value = $aIdxVals
$aIdxRng.Font.Bold = $true
$aIdxRng.Borders.LineStyle = 1
$aIdxRng.HorizontalAlignment = -4108
$aIdxRng.VerticalAlignment = -4160

$summary.Range("A1").Select()
